$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    1.8214061808877799,
    1.8227503596385899,
    1.8265227661171199,
    1.8309108405068999,
    1.83453322565197,
    1.83640958041136,
    1.84655529965402,
    1.86178787537039,
    1.8761829819193301,
    1.8882659536581401,
    1.9021695443261899,
    1.9144027586515,
    1.9264952666288699,
    1.9364271487228799,
    1.9417890360740699,
    1.93927316537913,
    1.93445753230775,
    1.9307922071905901,
    1.92638356351583,
    1.92037486267074,
    1.9139591929015001,
    1.9079505249697799,
    1.90339349685881,
    1.89934579822159,
    1.8950939885446001,
    1.8905862192717,
    1.88558090471816,
    1.8825240164438599,
    1.8806564194393001,
    1.8788605685061901,
    1.8757283731324501,
    1.87350918711953,
    1.8695727601415999,
    1.8637865042948001,
    1.8562143764576799,
    1.8464443455908099,
    1.84033542429318,
    1.83738597184403,
    1.8357571400372601,
    1.83422465878323,
    1.83433767392547,
    1.8336877625638,
    1.83064314451218,
    1.82675074679661,
    1.82333879196735,
    1.81887891549641,
    1.8139701128072201,
    1.8112887812363201,
    1.8099328599834399,
    1.8086509408379601,
    1.8078796362334899,
    1.8072999104527401,
    1.8070419808063201,
    1.8069674568077501,
    1.80688997060288,
    1.8067440306418401,
    1.8060207369756101,
    1.8071107271678599,
    1.8097614415729999,
    1.81313476638433,
    1.8175976262329501,
    1.8206970772481601,
    1.82495121889568,
    1.8317225324875801,
    1.8414378213860301,
    1.85043722923334,
    1.85859657911924,
    1.8655131427307801,
    1.8708716549816999,
    1.8747302345538599,
    1.8798133800441901,
    1.8828343045485001,
    1.88575248231366,
    1.8895666013553201,
    1.89471757904767,
    1.8992592311940699,
    1.9073888709620801,
    1.91597336617592,
    1.92565122343256,
    1.93846940387767,
    1.9499921790306001,
    1.9615763056829001,
    1.9703078541386501,
    1.97584563960826,
    1.9793426608873701,
    1.98002345245355,
    1.9826512892286401,
    1.9873373378077299,
    1.9942170102322101,
    2.00323343715195,
    2.0096527125913699,
    2.0108075523287798,
    2.0057871913184702,
    1.9957377288935401,
    1.98357210908476,
    1.9669521605540099,
    1.94853776379348,
    1.9302028993814599,
    1.9124996149001801,
    1.8947963304188999
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $values[$i]
}
